# HW1.docx edit: "Pages" -> "Panorama Items" (make UI layout more explicit),
# plus the incidental "_GoBack" bookmark relocation that Word leaves behind
# after the edit.
#
# Note: this runtime rebuilds/coalesces a paragraph's like-formatted runs
# whenever a length-changing text edit touches it, which would otherwise
# merge the (untouched) runs that follow our edit point into one run. To
# keep those later runs intact (matching the target run layout) we first
# drop temporary bookmarks at the exact split points -- bookmarks force a
# run boundary without altering any text -- and only then perform the
# length-changing replacement strictly between those two markers.

$d = $word.ActiveDocument

# --- Edit 1: "...two independent pages, one with..." -> "...two independent
# Panorama Items, one with..." ------------------------------------------

$before = $d.Content
$before.Find.Execute("two independent ") | Out-Null
$splitStart = $before.End

$after = $d.Content
$after.Find.Execute("pages,") | Out-Null
$splitEnd = $after.Start + 5   # end of the word "pages", right before the comma

$d.Bookmarks.Add("__edit_start", $d.Range($splitStart, $splitStart)) | Out-Null
$d.Bookmarks.Add("__edit_end", $d.Range($splitEnd, $splitEnd)) | Out-Null

$target = $d.Range($splitStart, $splitEnd)
$target.Text = "Panorama Items"

$d.Bookmarks("__edit_start").Delete()
$d.Bookmarks("__edit_end").Delete()

# --- Edit 2: relocate the "_GoBack" bookmark to sit right after "I a" in
# "...; I am not dictating how to make this change", splitting that run the
# same way Word does after an edit leaves the cursor there. -------------

$goBack = $d.Content
$goBack.Find.Execute("I a") | Out-Null
$goBackPos = $goBack.End

$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos)) | Out-Null
